$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 1.139906333333333
$ws.Cells.Item(2, 8).Value = 3.419719
$ws.Cells.Item(2, 9).Value = 0.2178538649973528
$ws.Cells.Item(2, 10).Value = 0.2178538649973527
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 35.48871133333333
$ws.Cells.Item(2, 14).Value = 106.466134
$ws.Cells.Item(2, 15).Value = 0.7670904531193554
$ws.Cells.Item(2, 16).Value = 0.7670904531193554
$ws.Cells.Item(2, 17).Value = 40.45380681070511
$ws.Cells.Item(2, 18).Value = 364.084261296346
$ws.Cells.Item(2, 19).Value = 0.1671136200146222
$ws.Cells.Item(2, 20).Value = 0.1671136200146222

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 1.139906333333333
$ws.Cells.Item(3, 8).Value = 3.419719
$ws.Cells.Item(3, 9).Value = 0.2178538649973528
$ws.Cells.Item(3, 10).Value = 0.2178538649973527
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 2.613261333333333
$ws.Cells.Item(3, 14).Value = 7.839784
$ws.Cells.Item(3, 15).Value = 0.05648578787427251
$ws.Cells.Item(3, 16).Value = 0.0564857878742725
$ws.Cells.Item(3, 17).Value = 2.978873144521778
$ws.Cells.Item(3, 18).Value = 26.809858300696
$ws.Cells.Item(3, 19).Value = 0.01230564720583087
$ws.Cells.Item(3, 20).Value = 0.01230564720583086

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 1.139906333333333
$ws.Cells.Item(4, 8).Value = 3.419719
$ws.Cells.Item(4, 9).Value = 0.2178538649973528
$ws.Cells.Item(4, 10).Value = 0.2178538649973527
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.596082
$ws.Cells.Item(4, 14).Value = 1.788246
$ws.Cells.Item(4, 15).Value = 0.01288434531142903
$ws.Cells.Item(4, 16).Value = 0.01288434531142903
$ws.Cells.Item(4, 17).Value = 0.679477646986
$ws.Cells.Item(4, 18).Value = 6.115298822873999
$ws.Cells.Item(4, 19).Value = 0.002806904424055335
$ws.Cells.Item(4, 20).Value = 0.002806904424055334

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 1.139906333333333
$ws.Cells.Item(5, 8).Value = 3.419719
$ws.Cells.Item(5, 9).Value = 0.2178538649973528
$ws.Cells.Item(5, 10).Value = 0.2178538649973527
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.482776
$ws.Cells.Item(5, 14).Value = 1.448328
$ws.Cells.Item(5, 15).Value = 0.01043522987117622
$ws.Cells.Item(5, 16).Value = 0.01043522987117622
$ws.Cells.Item(5, 17).Value = 0.5503194199813334
$ws.Cells.Item(5, 18).Value = 4.952874779831999
$ws.Cells.Item(5, 19).Value = 0.002273355159571566
$ws.Cells.Item(5, 20).Value = 0.002273355159571566

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1.139906333333333
$ws.Cells.Item(6, 8).Value = 3.419719
$ws.Cells.Item(6, 9).Value = 0.2178538649973528
$ws.Cells.Item(6, 10).Value = 0.2178538649973527
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 7.083219666666667
$ws.Cells.Item(6, 14).Value = 21.249659
$ws.Cells.Item(6, 15).Value = 0.1531041838237668
$ws.Cells.Item(6, 16).Value = 0.1531041838237668
$ws.Cells.Item(6, 17).Value = 8.074206958424556
$ws.Cells.Item(6, 18).Value = 72.667862625821
$ws.Cells.Item(6, 19).Value = 0.03335433819327278
$ws.Cells.Item(6, 20).Value = 0.03335433819327276

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1.902924
$ws.Cells.Item(7, 8).Value = 5.708772
$ws.Cells.Item(7, 9).Value = 0.3636784322304457
$ws.Cells.Item(7, 10).Value = 0.3636784322304456
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 35.48871133333333
$ws.Cells.Item(7, 14).Value = 106.466134
$ws.Cells.Item(7, 15).Value = 0.7670904531193554
$ws.Cells.Item(7, 16).Value = 0.7670904531193554
$ws.Cells.Item(7, 17).Value = 67.532320525272
$ws.Cells.Item(7, 18).Value = 607.7908847274481
$ws.Cells.Item(7, 19).Value = 0.2789742533693894
$ws.Cells.Item(7, 20).Value = 0.2789742533693893

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1.902924
$ws.Cells.Item(8, 8).Value = 5.708772
$ws.Cells.Item(8, 9).Value = 0.3636784322304457
$ws.Cells.Item(8, 10).Value = 0.3636784322304456
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 2.613261333333333
$ws.Cells.Item(8, 14).Value = 7.839784
$ws.Cells.Item(8, 15).Value = 0.05648578787427251
$ws.Cells.Item(8, 16).Value = 0.0564857878742725
$ws.Cells.Item(8, 17).Value = 4.972837709472
$ws.Cells.Item(8, 18).Value = 44.755539385248
$ws.Cells.Item(8, 19).Value = 0.02054266277741694
$ws.Cells.Item(8, 20).Value = 0.02054266277741694

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1.902924
$ws.Cells.Item(9, 8).Value = 5.708772
$ws.Cells.Item(9, 9).Value = 0.3636784322304457
$ws.Cells.Item(9, 10).Value = 0.3636784322304456
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.596082
$ws.Cells.Item(9, 14).Value = 1.788246
$ws.Cells.Item(9, 15).Value = 0.01288434531142903
$ws.Cells.Item(9, 16).Value = 0.01288434531142903
$ws.Cells.Item(9, 17).Value = 1.134298743768
$ws.Cells.Item(9, 18).Value = 10.208688693912
$ws.Cells.Item(9, 19).Value = 0.004685758503176202
$ws.Cells.Item(9, 20).Value = 0.004685758503176201

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 1.902924
$ws.Cells.Item(10, 8).Value = 5.708772
$ws.Cells.Item(10, 9).Value = 0.3636784322304457
$ws.Cells.Item(10, 10).Value = 0.3636784322304456
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 0.482776
$ws.Cells.Item(10, 14).Value = 1.448328
$ws.Cells.Item(10, 15).Value = 0.01043522987117622
$ws.Cells.Item(10, 16).Value = 0.01043522987117622
$ws.Cells.Item(10, 17).Value = 0.918686037024
$ws.Cells.Item(10, 18).Value = 8.268174333216
$ws.Cells.Item(10, 19).Value = 0.003795068039513682
$ws.Cells.Item(10, 20).Value = 0.003795068039513681

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 1.902924
$ws.Cells.Item(11, 8).Value = 5.708772
$ws.Cells.Item(11, 9).Value = 0.3636784322304457
$ws.Cells.Item(11, 10).Value = 0.3636784322304456
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 7.083219666666667
$ws.Cells.Item(11, 14).Value = 21.249659
$ws.Cells.Item(11, 15).Value = 0.1531041838237668
$ws.Cells.Item(11, 16).Value = 0.1531041838237668
$ws.Cells.Item(11, 17).Value = 13.478828700972
$ws.Cells.Item(11, 18).Value = 121.309458308748
$ws.Cells.Item(11, 19).Value = 0.05568068954094947
$ws.Cells.Item(11, 20).Value = 0.05568068954094945

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 2.189605333333333
$ws.Cells.Item(12, 8).Value = 6.568816
$ws.Cells.Item(12, 9).Value = 0.4184677027722017
$ws.Cells.Item(12, 10).Value = 0.4184677027722016
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 35.48871133333333
$ws.Cells.Item(12, 14).Value = 106.466134
$ws.Cells.Item(12, 15).Value = 0.7670904531193554
$ws.Cells.Item(12, 16).Value = 0.7670904531193554
$ws.Cells.Item(12, 17).Value = 77.70627160859378
$ws.Cells.Item(12, 18).Value = 699.3564444773441
$ws.Cells.Item(12, 19).Value = 0.321002579735344
$ws.Cells.Item(12, 20).Value = 0.3210025797353439

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 2.189605333333333
$ws.Cells.Item(13, 8).Value = 6.568816
$ws.Cells.Item(13, 9).Value = 0.4184677027722017
$ws.Cells.Item(13, 10).Value = 0.4184677027722016
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 2.613261333333333
$ws.Cells.Item(13, 14).Value = 7.839784
$ws.Cells.Item(13, 15).Value = 0.05648578787427251
$ws.Cells.Item(13, 16).Value = 0.0564857878742725
$ws.Cells.Item(13, 17).Value = 5.722010952860444
$ws.Cells.Item(13, 18).Value = 51.498098575744
$ws.Cells.Item(13, 19).Value = 0.0236374778910247
$ws.Cells.Item(13, 20).Value = 0.0236374778910247

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 2.189605333333333
$ws.Cells.Item(14, 8).Value = 6.568816
$ws.Cells.Item(14, 9).Value = 0.4184677027722017
$ws.Cells.Item(14, 10).Value = 0.4184677027722016
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 0.596082
$ws.Cells.Item(14, 14).Value = 1.788246
$ws.Cells.Item(14, 15).Value = 0.01288434531142903
$ws.Cells.Item(14, 16).Value = 0.01288434531142903
$ws.Cells.Item(14, 17).Value = 1.305184326304
$ws.Cells.Item(14, 18).Value = 11.746658936736
$ws.Cells.Item(14, 19).Value = 0.005391682384197493
$ws.Cells.Item(14, 20).Value = 0.005391682384197492

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 2.189605333333333
$ws.Cells.Item(15, 8).Value = 6.568816
$ws.Cells.Item(15, 9).Value = 0.4184677027722017
$ws.Cells.Item(15, 10).Value = 0.4184677027722016
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.482776
$ws.Cells.Item(15, 14).Value = 1.448328
$ws.Cells.Item(15, 15).Value = 0.01043522987117622
$ws.Cells.Item(15, 16).Value = 0.01043522987117622
$ws.Cells.Item(15, 17).Value = 1.057088904405333
$ws.Cells.Item(15, 18).Value = 9.513800139648
$ws.Cells.Item(15, 19).Value = 0.00436680667209097
$ws.Cells.Item(15, 20).Value = 0.004366806672090969

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 2.189605333333333
$ws.Cells.Item(16, 8).Value = 6.568816
$ws.Cells.Item(16, 9).Value = 0.4184677027722017
$ws.Cells.Item(16, 10).Value = 0.4184677027722016
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 7.083219666666667
$ws.Cells.Item(16, 14).Value = 21.249659
$ws.Cells.Item(16, 15).Value = 0.1531041838237668
$ws.Cells.Item(16, 16).Value = 0.1531041838237668
$ws.Cells.Item(16, 17).Value = 15.50945555930489
$ws.Cells.Item(16, 18).Value = 139.585100033744
$ws.Cells.Item(16, 19).Value = 0.06406915608954458
$ws.Cells.Item(16, 20).Value = 0.06406915608954455
